# Fruta / hortaliza, semanal
# Update existing rows 2-13 with refreshed weekly prices, and append a new
# row 14 (new market-price observation) to the Arándano (blue) sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 44594
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 2500
$ws.Range("O2").Value = 2800
$ws.Range("P2").Value = 2650
$ws.Range("S2").Value = 1325

# --- Row 3 ---
$ws.Range("D3").Value = 44539
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 3800
$ws.Range("O3").Value = 4000
$ws.Range("P3").Value = 3900
$ws.Range("Q3").Value = "$/bandeja 2 kilos"
$ws.Range("R3").Value = "Región del Maule"
$ws.Range("S3").Value = 1950
$ws.Range("T3").Value = 2

# --- Row 4 ---
$ws.Range("D4").Value = 44540
$ws.Range("M4").Value = 240
$ws.Range("N4").Value = 3500
$ws.Range("O4").Value = 3800
$ws.Range("P4").Value = 3650
$ws.Range("S4").Value = 1825

# --- Row 7 ---
$ws.Range("D7").Value = 44187
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 2800
$ws.Range("O7").Value = 3000
$ws.Range("P7").Value = 2900
$ws.Range("S7").Value = 1450

# --- Row 8 ---
$ws.Range("D8").Value = 44187
$ws.Range("M8").Value = 65
$ws.Range("N8").Value = 1400
$ws.Range("O8").Value = 1500
$ws.Range("P8").Value = 1446
$ws.Range("Q8").Value = "$/envase 1 kilo"
$ws.Range("R8").Value = "Provincia de Diguillín"
$ws.Range("S8").Value = 1446
$ws.Range("T8").Value = 1

# --- Row 9 ---
$ws.Range("D9").Value = 44174
$ws.Range("M9").Value = 150
$ws.Range("N9").Value = 3700
$ws.Range("O9").Value = 3800
$ws.Range("P9").Value = 3747
$ws.Range("R9").Value = "Provincia de Linares"
$ws.Range("S9").Value = 1874

# --- Row 10 ---
$ws.Range("D10").Value = 44937
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 2500
$ws.Range("O10").Value = 3000
$ws.Range("P10").Value = 2750
$ws.Range("R10").Value = "Provincia de Diguillín"
$ws.Range("S10").Value = 1375

# --- Row 11 ---
$ws.Range("D11").Value = 44932
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 3000
$ws.Range("O11").Value = 3000
$ws.Range("P11").Value = 3000
$ws.Range("S11").Value = 1500

# --- Row 12 ---
$ws.Range("M12").Value = 65
$ws.Range("N12").Value = 3600
$ws.Range("O12").Value = 3800
$ws.Range("P12").Value = 3692
$ws.Range("Q12").Value = "$/bandeja 2 kilos"
$ws.Range("S12").Value = 1846
$ws.Range("T12").Value = 2

# --- Row 13 ---
$ws.Range("D13").Value = 44181
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 1800
$ws.Range("O13").Value = 2000
$ws.Range("P13").Value = 1875
$ws.Range("Q13").Value = "$/envase 1 kilo"
$ws.Range("S13").Value = 1875
$ws.Range("T13").Value = 1

# --- New row 14 ---
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C14").Value = "Ñuble"
$ws.Range("D14").Value = 44931
$ws.Range("D14").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E14").Value = 16
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100101
$ws.Range("H14").Value = "Berries"
$ws.Range("I14").Value = 100101001
$ws.Range("J14").Value = "Arándano (blue)"
$ws.Range("K14").Value = "Sin especificar"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 3000
$ws.Range("O14").Value = 3000
$ws.Range("P14").Value = 3000
$ws.Range("Q14").Value = "$/bandeja 2 kilos"
$ws.Range("R14").Value = "Provincia de Diguillín"
$ws.Range("S14").Value = 1500
$ws.Range("T14").Value = 2
